$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.092.75"
$ws.Range("E2").Value = "  +6.34%  "
$ws.Range("D3").Value = "1.893.03"
$ws.Range("E3").Value = "  +5.89%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "249.71"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.5006"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("D8").Value = "45.87"
$ws.Range("E8").Value = "  +8.68%  "
$ws.Range("D9").Value = "0.2877"
$ws.Range("E9").Value = "  +6.86%  "
$ws.Range("D10").Value = "0.06571"
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("D11").Value = "1.886.38"
$ws.Range("E11").Value = "  +5.54%  "
$ws.Range("D12").Value = "17.27"
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("D13").Value = "0.07242"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").Value = "0.6705"
$ws.Range("E14").Value = "  +6.48%  "
$ws.Range("D15").Value = "85.10"
$ws.Range("E15").Value = "  +6.09%  "
$ws.Range("D16").Value = "4.836"
$ws.Range("E16").Value = "  +3.62%  "
$ws.Range("D17").Value = "30.114.84"
$ws.Range("E17").Value = "  +6.58%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  +7.05%  "
$ws.Range("D20").Value = "0.000007552"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +5.64%  "
$ws.Range("D23").Value = "4.786"
$ws.Range("E23").Value = "  +4.77%  "
$ws.Range("D24").Value = "5.556"
$ws.Range("E24").Value = "  +5.76%  "
$ws.Range("D25").Value = "9.059"
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("D26").Value = "145.39"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("D27").Value = "136.33"
$ws.Range("E27").Value = "  +24.01%  "
$ws.Range("D28").Value = "16.81"
$ws.Range("E28").Value = "  +6.43%  "
$ws.Range("D29").Value = "1.953"
$ws.Range("E29").Value = "  +5.32%  "
$ws.Range("D30").Value = "1.373"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "4.207"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").Value = "0.08673"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("D33").Value = "3.942"
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("D34").Value = "0.05036"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("D36").Value = "0.6952"
$ws.Range("E36").Value = "  +5.46%  "
$ws.Range("D37").Value = "2.689"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("D38").Value = "2.285"
$ws.Range("E38").Value = "  +10.49%  "
$ws.Range("D39").Value = "2.778"
$ws.Range("E39").Value = "  +6.29%  "
$ws.Range("D40").Value = "0.9632"
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("D41").Value = "0.01645"
$ws.Range("E41").Value = "  +5.85%  "
$ws.Range("D42").Value = "6.065"
$ws.Range("E42").Value = "  +2.68%  "
$ws.Range("D43").Value = "105.75"
$ws.Range("E43").Value = "  +5.93%  "
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "0.4226"
$ws.Range("E45").Value = "  +5.66%  "
$ws.Range("D46").Value = "7.464"
$ws.Range("E46").Value = "  +3.66%  "
$ws.Range("D47").Value = "0.1258"
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("D48").Value = "0.05667"
$ws.Range("D49").Value = "32.60"
$ws.Range("E49").Value = "  +5.91%  "
$ws.Range("D50").Value = "8.299"
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("D51").Value = "0.3735"
$ws.Range("E51").Value = "  +6.72%  "
